# Updated export from Jun's tool
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet: refresh the source citation block and restore the Note text
# to reference 1-year (instead of 3-year) elasticities.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("B4").Value = 2005
$about.Range("B5").Value = "Price Responsiveness in the AEO 2003 NEMS Residential and Commercial Buildings Sector Models"

# Swap the hyperlink target + its displayed text (drop the old relationship first
# so we don't leave a stale one behind).
$about.Range("B6").Hyperlinks.Delete()
$about.Hyperlinks.Add($about.Range("B6"), "http://www.eia.gov/oiaf/analysispaper/elasticity/pdf/tbl.pdf")
$about.Range("B6").Value = "http://www.eia.gov/oiaf/analysispaper/elasticity/pdf/tbl.pdf"

$about.Range("B7").Value = "Page 1, Table 1"

$about.Range("A10").Value = "We use same-price, long-run elasticities minus the 1-year short-run elasticities."
$about.Range("A11").Value = "We calculate it this way because we assume that 1-year elasticities primarily reflect behavior"
$about.Range("A14").Value = "all timescales.  So, the portion of the long-run elasticitiy represented by the 1-year elasticity"

# ---------------------------------------------------------------------------
# "EIA Table 1" sheet: refreshed elasticity data export.
# ---------------------------------------------------------------------------
$eia = $wb.Worksheets.Item("EIA Table 1")

# Residential block (rows 7-9)
$eia.Range("B7").Value = -0.2
$eia.Range("C7").Value = -0.28999999999999998
$eia.Range("D7").Value = -0.34
$eia.Range("E7").Value = -0.49
$eia.Range("F7").Value = 0.01
$eia.Range("G7").Value = 0

$eia.Range("B8").Value = -0.14000000000000001
$eia.Range("C8").Value = -0.24
$eia.Range("D8").Value = -0.3
$eia.Range("E8").Value = 0.13
$eia.Range("F8").Value = -0.41
$eia.Range("G8").Value = 0.02

$eia.Range("B9").Value = -0.15
$eia.Range("C9").Value = -0.27
$eia.Range("D9").Value = -0.34
$eia.Range("E9").Value = 0.01
$eia.Range("F9").Value = 0.05
$eia.Range("G9").Value = -0.6

# Commercial block (rows 14-16)
$eia.Range("B14").Value = -0.1
$eia.Range("C14").Value = -0.17
$eia.Range("D14").Value = -0.2
$eia.Range("E14").Value = -0.45
$eia.Range("F14").Value = 0.01
$eia.Range("G14").Value = 0

$eia.Range("B15").Value = -0.14000000000000001
$eia.Range("C15").Value = -0.24
$eia.Range("D15").Value = -0.28999999999999998
$eia.Range("E15").Value = 0.86
$eia.Range("F15").Value = -0.4
$eia.Range("G15").Value = 0.01

$eia.Range("B16").Value = -0.13
$eia.Range("C16").Value = -0.23
$eia.Range("D16").Value = -0.28000000000000003
$eia.Range("E16").Value = 0.08
$eia.Range("F16").Value = 0.75
$eia.Range("G16").Value = -0.39

# ---------------------------------------------------------------------------
# "EoCEDwEC" sheet: the short-run elasticity subtraction now uses the 1-year
# column (B) instead of the 3-year column (D).
# ---------------------------------------------------------------------------
$eoc = $wb.Worksheets.Item("EoCEDwEC")

$eoc.Range("B2").Formula = "='EIA Table 1'!E7-'EIA Table 1'!B7"
$eoc.Range("D2").Formula = "='EIA Table 1'!E14-'EIA Table 1'!B14"
$eoc.Range("B4").Formula = "='EIA Table 1'!F8-'EIA Table 1'!B8"
$eoc.Range("D4").Formula = "='EIA Table 1'!F15-'EIA Table 1'!B15"
$eoc.Range("B5").Formula = "='EIA Table 1'!G9-'EIA Table 1'!B9"
$eoc.Range("D5").Formula = "='EIA Table 1'!G16-'EIA Table 1'!B16"

# ---------------------------------------------------------------------------
# View state: active tab moves from "EoCEDwEC" (index 2) to "About" (index 0).
# ---------------------------------------------------------------------------
$about.Activate()

$wb.Application.Calculate()
